$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) switch from the deck's custom table
#    style to the built-in "Medium Style 2 - Accent 1"-family style whose
#    GUID is {BAFECB4A-1DB8-4106-8E11-C5FE053E26C7}.
# ---------------------------------------------------------------------------
$newTableStyle = "{BAFECB4A-1DB8-4106-8E11-C5FE053E26C7}"
$tableSlideNumbers = @(14, 15, 16)

foreach ($slideNum in $tableSlideNumbers) {
    $slide = $p.Slides.Item($slideNum)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme (the "Integral" / Red Violet colour scheme used by
#    the slide master) is swapped for the plain built-in "Office Theme"
#    colour scheme.
# ---------------------------------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$firstSlide = $p.Slides.Item(1)
$themeColors = $firstSlide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
